$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "29.629.10"
$ws.Range("E2").Value = "  -0.92%  "

$ws.Range("D3").Value = "1.860.66"
$ws.Range("E3").Value = "  -0.08%  "

$ws.Range("D4").Value = "0.9982"
$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").Value = "242.66"
$ws.Range("E5").Value = "  -1.04%  "

$ws.Range("D6").Value = "0.6356"
$ws.Range("E6").Value = "  -3.70%  "

$ws.Range("D7").Value = "0.9992"
$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").Value = "0.07617"
$ws.Range("E8").Value = "  +0.49%  "

$ws.Range("D9").Value = "0.3003"
$ws.Range("E9").Value = "  +0.93%  "

$ws.Range("D10").Value = "24.72"
$ws.Range("E10").Value = "  -0.19%  "

$ws.Range("D11").Value = "0.07742"
$ws.Range("E11").Value = "  +1.02%  "

$ws.Range("D12").Value = "1.911.55"
$ws.Range("E12").Value = "  +2.41%  "

$ws.Range("D13").Value = "0.6974"
$ws.Range("E13").Value = "  +0.56%  "

$ws.Range("D14").Value = "5.042"
$ws.Range("E14").Value = "  -0.70%  "

$ws.Range("D15").Value = "83.94"
$ws.Range("E15").Value = "  -0.04%  "

$ws.Range("D16").Value = "0.00001000"
$ws.Range("E16").Value = "  +3.04%  "

$ws.Range("D17").Value = "2.131.64"
$ws.Range("E17").Value = "  +0.35%  "

$ws.Range("D18").Value = "6.281"
$ws.Range("E18").Value = "  +2.29%  "

$ws.Range("D19").Value = "29.608.47"
$ws.Range("E19").Value = "  -1.06%  "

$ws.Range("D20").Value = "236.08"
$ws.Range("E20").Value = "  -0.26%  "

$ws.Range("D21").Value = "12.62"
$ws.Range("E21").Value = "  -0.85%  "

$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").Value = "7.689"
$ws.Range("E23").Value = "  -0.78%  "

$ws.Range("D24").Value = "0.9993"
$ws.Range("E24").Value = "  -0.29%  "

$ws.Range("D25").Value = "156.29"
$ws.Range("E25").Value = "  -1.49%  "

$ws.Range("D26").Value = "0.1405"
$ws.Range("E26").Value = "  -3.11%  "

$ws.Range("D27").Value = "8.512"
$ws.Range("E27").Value = "  -1.09%  "

$ws.Range("E28").Value = "  -0.97%  "

$ws.Range("D29").Value = "1.479"
$ws.Range("E29").Value = "  -1.24%  "

$ws.Range("D30").Value = "0.05821"
$ws.Range("E30").Value = "  -4.06%  "

$ws.Range("E31").Value = "  -1.88%  "

$ws.Range("D32").Value = "4.142"
$ws.Range("E32").Value = "  -0.72%  "

$ws.Range("D33").Value = "4.047"
$ws.Range("E33").Value = "  -1.48%  "

$ws.Range("D34").Value = "1.913"
$ws.Range("E34").Value = "  +1.57%  "

$ws.Range("D35").Value = "1.174"
$ws.Range("E35").Value = "  -0.82%  "

$ws.Range("D36").Value = "0.7269"
$ws.Range("E36").Value = "  -1.42%  "

$ws.Range("D37").Value = "2.587"
$ws.Range("E37").Value = "  -1.04%  "

$ws.Range("D38").Value = "1.259.59"
$ws.Range("E38").Value = "  +3.87%  "

$ws.Range("D39").Value = "2.813"
$ws.Range("E39").Value = "  -0.23%  "

$ws.Range("D40").Value = "0.01813"
$ws.Range("E40").Value = "  +0.75%  "

$ws.Range("D41").Value = "0.9088"
$ws.Range("E41").Value = "  -0.83%  "

$ws.Range("D42").Value = "6.155"
$ws.Range("E42").Value = "  -3.12%  "

$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").Value = "2.029.98"
$ws.Range("E43").Value = "  -0.24%  "

$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "0.9992"
$ws.Range("E44").Value = "  -0.17%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "68.49"
$ws.Range("E45").Value = "  +1.28%  "

$ws.Range("D46").Value = "101.67"
$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("D47").Value = "7.380"
$ws.Range("E47").Value = "  -5.78%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "9.255"
$ws.Range("E48").Value = "  +0.49%  "

$ws.Range("D49").Value = "0.4070"
$ws.Range("E49").Value = "  -0.42%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.00000000118"
$ws.Range("E50").Value = "  -3.77%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "1.722"
$ws.Range("E51").Value = "  +2.17%  "
